# Auto-generated Excel COM-interop script
# Adds columns J:O (SHG-CGG, SHG-Jumbo, Lead-Min 99.99%, EPG, HZDA3, HZDA5),
# fixes a handful of cells that should hold numeric 0 instead of text "0",
# and appends 14 new data rows (104-117) sourced from the 2023-03-20 circular.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New header cells J1:O1, matching the style of the existing header row ---
$ws.Range("I1").Copy() | Out-Null
$ws.Range("J1:O1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("J1").Value = "SHG-`nCGG"
$ws.Range("K1").Value = "SHG-`nJumbo"
$ws.Range("L1").Value = "Lead`n-Min`n99.99`n%"
$ws.Range("M1").Value = "EPG"
$ws.Range("N1").Value = "HZDA3`n(Zn`nAlloy`n)"
$ws.Range("O1").Value = "HZDA5`n(Zn`nAlloy`n)"

# --- 2. A few pre-existing "0" text cells turn out to be numeric zeros ---
$ws.Range("H90").Value = 0
$ws.Range("H91").Value = 0
$ws.Range("F92").Value = 0
$ws.Range("H92").Value = 0

# --- 3. Append rows 104-117 (2023-03-20 HZL circular) across columns A:O ---
# Row 104
$ws.Range("A104").Value = "'2023-03-20"
$ws.Range("B104").Value = "https://rashtriyametal.com/wp-content/uploads/2023/03/HZL20032023.pdf"
$ws.Range("C104").Value = "CHANDERIYA`nLEAD ZINC`nSMELTER"
$ws.Range("D104").Value = "'281,000"
$ws.Range("G104").Value = "'280,500"
$ws.Range("H104").Value = "'279,000"
$ws.Range("J104").Value = "'282,800"
$ws.Range("K104").Value = "'281,500"
$ws.Range("L104").Value = "'203,900"
$ws.Range("M104").Value = "'282,500"
$ws.Range("N104").Value = "'284,000"
$ws.Range("O104").Value = "'287,000"

# Row 105
$ws.Range("A105").Value = "'2023-03-20"
$ws.Range("B105").Value = "https://rashtriyametal.com/wp-content/uploads/2023/03/HZL20032023.pdf"
$ws.Range("C105").Value = "HYDRO-1 UNIT"
$ws.Range("D105").Value = "'281,000"
$ws.Range("G105").Value = "'280,500"
$ws.Range("H105").Value = "'279,000"
$ws.Range("J105").Value = "'282,800"
$ws.Range("K105").Value = "'281,500"
$ws.Range("L105").Value = "'203,900"
$ws.Range("M105").Value = "'282,500"
$ws.Range("N105").Value = "'284,000"
$ws.Range("O105").Value = "'287,000"

# Row 106
$ws.Range("A106").Value = "'2023-03-20"
$ws.Range("B106").Value = "https://rashtriyametal.com/wp-content/uploads/2023/03/HZL20032023.pdf"
$ws.Range("C106").Value = "NEW HYDRO`nSMELTER`nCHANDERIYA"
$ws.Range("D106").Value = "'281,000"
$ws.Range("G106").Value = "'280,500"
$ws.Range("H106").Value = "'279,000"
$ws.Range("J106").Value = "'282,800"
$ws.Range("K106").Value = "'281,500"
$ws.Range("L106").Value = "'203,900"
$ws.Range("M106").Value = "'282,500"
$ws.Range("N106").Value = "'284,000"
$ws.Range("O106").Value = "'287,000"

# Row 107
$ws.Range("A107").Value = "'2023-03-20"
$ws.Range("B107").Value = "https://rashtriyametal.com/wp-content/uploads/2023/03/HZL20032023.pdf"
$ws.Range("C107").Value = "ZINC SMELTER`nDEBRI"
$ws.Range("D107").Value = "'0"
$ws.Range("G107").Value = "'280,500"
$ws.Range("H107").Value = "'0"
$ws.Range("J107").Value = "'0"
$ws.Range("K107").Value = "'0"
$ws.Range("L107").Value = "'0"
$ws.Range("M107").Value = "'0"
$ws.Range("N107").Value = "'0"
$ws.Range("O107").Value = "'0"

# Row 108
$ws.Range("A108").Value = "'2023-03-20"
$ws.Range("B108").Value = "https://rashtriyametal.com/wp-content/uploads/2023/03/HZL20032023.pdf"
$ws.Range("C108").Value = "Pantnagar`nMelting&Castin`ngPlant"
$ws.Range("D108").Value = "'281,500"
$ws.Range("G108").Value = "'281,000"
$ws.Range("H108").Value = "'0"
$ws.Range("J108").Value = "'283,300"
$ws.Range("K108").Value = "'282,000"
$ws.Range("L108").Value = "'203,900"
$ws.Range("M108").Value = "'283,000"
$ws.Range("N108").Value = "'284,500"
$ws.Range("O108").Value = "'287,500"

# Row 109
$ws.Range("A109").Value = "'2023-03-20"
$ws.Range("B109").Value = "https://rashtriyametal.com/wp-content/uploads/2023/03/HZL20032023.pdf"
$ws.Range("C109").Value = "RAJPURA DARIBA`nLEAD SMELTER"
$ws.Range("D109").Value = "'0"
$ws.Range("G109").Value = "'0"
$ws.Range("H109").Value = "'0"
$ws.Range("J109").Value = "'0"
$ws.Range("K109").Value = "'0"
$ws.Range("L109").Value = "'203,900"
$ws.Range("M109").Value = "'0"
$ws.Range("N109").Value = "'0"
$ws.Range("O109").Value = "'0"

# Row 110
$ws.Range("A110").Value = "'2023-03-20"
$ws.Range("B110").Value = "https://rashtriyametal.com/wp-content/uploads/2023/03/HZL20032023.pdf"
$ws.Range("C110").Value = "Faridabad`nDepot"
$ws.Range("D110").Value = "'282,000"
$ws.Range("G110").Value = "'281,500"
$ws.Range("H110").Value = "'280,000"
$ws.Range("J110").Value = "'283,800"
$ws.Range("K110").Value = "'282,500"
$ws.Range("L110").Value = "'208,400"
$ws.Range("M110").Value = "'283,500"
$ws.Range("N110").Value = "'285,000"
$ws.Range("O110").Value = "'288,000"

# Row 111
$ws.Range("A111").Value = "'2023-03-20"
$ws.Range("B111").Value = "https://rashtriyametal.com/wp-content/uploads/2023/03/HZL20032023.pdf"
$ws.Range("C111").Value = "Panvel Depot"
$ws.Range("D111").Value = "'281,000"
$ws.Range("G111").Value = "'280,500"
$ws.Range("H111").Value = "'279,000"
$ws.Range("J111").Value = "'282,800"
$ws.Range("K111").Value = "'281,500"
$ws.Range("L111").Value = "'207,500"
$ws.Range("M111").Value = "'282,500"
$ws.Range("N111").Value = "'284,000"
$ws.Range("O111").Value = "'287,000"

# Row 112
$ws.Range("A112").Value = "'2023-03-20"
$ws.Range("B112").Value = "https://rashtriyametal.com/wp-content/uploads/2023/03/HZL20032023.pdf"
$ws.Range("C112").Value = "Baroda Depot"
$ws.Range("D112").Value = "'281,500"
$ws.Range("G112").Value = "'281,000"
$ws.Range("H112").Value = "'279,500"
$ws.Range("J112").Value = "'283,300"
$ws.Range("K112").Value = "'282,000"
$ws.Range("L112").Value = "'207,900"
$ws.Range("M112").Value = "'283,000"
$ws.Range("N112").Value = "'284,500"
$ws.Range("O112").Value = "'287,500"

# Row 113
$ws.Range("A113").Value = "'2023-03-20"
$ws.Range("B113").Value = "https://rashtriyametal.com/wp-content/uploads/2023/03/HZL20032023.pdf"
$ws.Range("C113").Value = "JAMSHEDPUR`nSTOCK POINT"
$ws.Range("D113").Value = "'281,000"
$ws.Range("G113").Value = "'280,500"
$ws.Range("H113").Value = "'279,000"
$ws.Range("J113").Value = "'282,800"
$ws.Range("K113").Value = "'281,500"
$ws.Range("L113").Value = "'206,900"
$ws.Range("M113").Value = "'282,500"
$ws.Range("N113").Value = "'284,000"
$ws.Range("O113").Value = "'287,000"

# Row 114
$ws.Range("A114").Value = "'2023-03-20"
$ws.Range("B114").Value = "https://rashtriyametal.com/wp-content/uploads/2023/03/HZL20032023.pdf"
$ws.Range("C114").Value = "Kolkata Depot"
$ws.Range("D114").Value = "'281,000"
$ws.Range("G114").Value = "'280,500"
$ws.Range("H114").Value = "'279,000"
$ws.Range("J114").Value = "'282,800"
$ws.Range("K114").Value = "'281,500"
$ws.Range("L114").Value = "'206,900"
$ws.Range("M114").Value = "'282,500"
$ws.Range("N114").Value = "'284,000"
$ws.Range("O114").Value = "'287,000"

# Row 115
$ws.Range("A115").Value = "'2023-03-20"
$ws.Range("B115").Value = "https://rashtriyametal.com/wp-content/uploads/2023/03/HZL20032023.pdf"
$ws.Range("C115").Value = "Bangalore`nDepot"
$ws.Range("D115").Value = "'281,500"
$ws.Range("G115").Value = "'281,000"
$ws.Range("H115").Value = "'279,500"
$ws.Range("J115").Value = "'283,300"
$ws.Range("K115").Value = "'282,000"
$ws.Range("L115").Value = "'206,400"
$ws.Range("M115").Value = "'283,000"
$ws.Range("N115").Value = "'284,500"
$ws.Range("O115").Value = "'287,500"

# Row 116
$ws.Range("A116").Value = "'2023-03-20"
$ws.Range("B116").Value = "https://rashtriyametal.com/wp-content/uploads/2023/03/HZL20032023.pdf"
$ws.Range("C116").Value = "Hyderabad`nDepot"
$ws.Range("D116").Value = "'281,500"
$ws.Range("G116").Value = "'281,000"
$ws.Range("H116").Value = "'279,500"
$ws.Range("J116").Value = "'283,300"
$ws.Range("K116").Value = "'282,000"
$ws.Range("L116").Value = "'207,900"
$ws.Range("M116").Value = "'283,000"
$ws.Range("N116").Value = "'284,500"
$ws.Range("O116").Value = "'287,500"

# Row 117
$ws.Range("A117").Value = "'2023-03-20"
$ws.Range("B117").Value = "https://rashtriyametal.com/wp-content/uploads/2023/03/HZL20032023.pdf"
$ws.Range("C117").Value = "Chennai Depot"
$ws.Range("D117").Value = "'281,000"
$ws.Range("G117").Value = "'280,500"
$ws.Range("H117").Value = "'279,000"
$ws.Range("J117").Value = "'282,800"
$ws.Range("K117").Value = "'281,500"
$ws.Range("L117").Value = "'208,400"
$ws.Range("M117").Value = "'282,500"
$ws.Range("N117").Value = "'284,000"
$ws.Range("O117").Value = "'287,000"

